# Activity.xlsx — add a handful of newly-logged dates to 工作表1 and
# update the sheet's view state (matches the "Add files via upload" commit).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("工作表1")

# --- new / updated date values -------------------------------------------------
# Row 4: start date for a new row item in column E (uses column E's existing format)
$ws.Range("E4").Value = 45931

# Row 5: fill in a previously-blank date in column R
$ws.Range("R5").Value = 45931

# Row 6: new date in column F
$ws.Range("F6").Value = 45924

# Row 7: fill in a previously-blank date in column O
$ws.Range("O7").Value = 45924

# Row 8: new date in column M, formatted like the rest of column M's existing
# "yyyy-mm-dd" cells (reuses the existing style used by sibling cells, e.g. M2/N2).
$ws.Range("M8").Value = 45931
$ws.Range("M8").NumberFormat = "yyyy\-mm\-dd;@"

# Row 9: new (still-empty) cell in column M that picks up its own short-date
# format — this mints a brand-new cell style (built-in date format 14) the
# same way Excel would after "Format Cells > Date > m/d/yyyy", while keeping
# the cell's existing fill.
$ws.Range("M9").NumberFormat = "mm-dd-yy"

# --- view state ------------------------------------------------------------
# Scroll the sheet so column F is the left-most visible column, then leave the
# selection on M8 (where the day's data entry finished).
$excel.ActiveWindow.ScrollColumn = 6
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("M8").Select()
